# Refresh the cryptos list: updated Price/Volume(1h) figures and the
# Chainlink / WrappedEther row swap (rows 15-16), per the Jan 15 2024 GitHub Action run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row updates: Row number -> column letter -> new text value.
# ("needsText" cells are pure numeric-looking strings like "1.00" or "6.70" that
#  Excel would otherwise silently coerce to a number, dropping the trailing zero /
#  formatting - force them to keep their literal text via the Text number format,
#  same as typing an apostrophe-prefixed value in the Excel UI.)
$updates = @(
    @{Row=2; D="42.752.59"; E="  +0.91%  "},
    @{Row=3; D="2.522.47"; E="  -0.02%  "},
    @{Row=4; D="1.00"; E="  +0.29%  "},
    @{Row=5; D="316.68"; E="  +4.10%  "},
    @{Row=6; D="94.57"; E="  -2.39%  "},
    @{Row=7; D="0.579"; E="  -1.10%  "},
    @{Row=8; E="  -0.02%  "},
    @{Row=9; D="0.529"; E="  -1.84%  "},
    @{Row=10; D="35.82"; E="  -2.11%  "},
    @{Row=11; D="0.0809"; E="  +0.03%  "},
    @{Row=12; D="7.55"; E="  -1.47%  "},
    @{Row=13; D="0.109"; E="  -2.84%  "},
    @{Row=14; D="2.912.12"; E="  +0.51%  "},
    @{Row=15; B="WrappedEther"; C="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D="2.529.41"; E="  -0.76%  "},
    @{Row=16; B="Chainlink"; C="https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D="15.19"; E="  +0.35%  "},
    @{Row=17; D="0.847"; E="  -1.81%  "},
    @{Row=18; D="42.857.42"; E="  +1.19%  "},
    @{Row=19; D="12.99"; E="  +0.07%  "},
    @{Row=20; D="6.70"; E="  +3.82%  "},
    @{Row=21; D="0.0₃0963"; E="  -0.91%  "},
    @{Row=22; D="69.73"; E="  -1.71%  "},
    @{Row=23; D="250.73"; E="  +0.06%  "},
    @{Row=24; E="  +1.12%  "},
    @{Row=25; E="  -0.05%  "},
    @{Row=26; D="26.72"; E="  -1.00%  "},
    @{Row=27; E="  -0.03%  "},
    @{Row=28; D="2.42"; E="  +4.01%  "},
    @{Row=29; D="40.27"; E="  +5.49%  "},
    @{Row=30; D="10.26"; E="  -0.07%  "},
    @{Row=31; D="5.97"; E="  +0.39%  "},
    @{Row=32; D="156.37"; E="  +0.80%  "},
    @{Row=33; E="  +2.13%  "},
    @{Row=34; D="18.96"; E="  +1.29%  "},
    @{Row=35; D="3.28"; E="  -1.42%  "},
    @{Row=36; D="0.0786"; E="  -0.35%  "},
    @{Row=37; E="  -0.42%  "},
    @{Row=38; E="  -2.13%  "},
    @{Row=39; E="  -0.55%  "},
    @{Row=40; D="23.58"; E="  -2.51%  "},
    @{Row=41; E="  +14.09%  "},
    @{Row=42; E="  +1.78%  "},
    @{Row=43; E="  +0.42%  "},
    @{Row=44; E="  -2.19%  "},
    @{Row=45; D="3.29"; E="  -2.99%  "},
    @{Row=46; D="2.019.38"; E="  -0.81%  "},
    @{Row=47; D="85.57"; E="  +1.07%  "},
    @{Row=48; D="8.76"; E="  -2.05%  "},
    @{Row=49; D="2.766.50"; E="  +0.23%  "},
    @{Row=50; D="73.53"; E="  +1.84%  "},
    @{Row=51; D="102.47"; E="  +0.55%  "}
)

foreach ($u in $updates) {
    foreach ($col in @("B", "C", "D", "E")) {
        if (-not $u.ContainsKey($col)) { continue }
        $value = $u[$col]
        $cell = $ws.Cells.Item($u.Row, @{B=2;C=3;D=4;E=5}[$col])
        if ($value -match "^[+-]?\d+(\.\d+)?$") {
            # Pure numeric-looking text - force Text format so Excel keeps it verbatim.
            $cell.NumberFormat = "@"
        }
        $cell.Value = $value
    }
}

Write-Output "Updated $($updates.Count) rows"
